$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New date serials (column A) for rows 2..48 after the upload/refresh.
$dates = @{
    2  = 45576
    3  = 45577
    4  = 45578
    5  = 45579
    6  = 45580
    7  = 45581
    8  = 45582
    9  = 45583
    10 = 45584
    11 = 45585
    12 = 45586
    13 = 45587
    14 = 45588
    15 = 45589
    16 = 45590
    17 = 45591
    18 = 45592
    19 = 45593
    20 = 45594
    21 = 45595
    22 = 45596
    23 = 45597
    24 = 45598
    25 = 45599
    26 = 45600
    27 = 45601
    28 = 45602
    29 = 45603
    30 = 45604
    31 = 45558
    32 = 45559
    33 = 45560
    34 = 45561
    35 = 45562
    36 = 45563
    37 = 45564
    38 = 45565
    39 = 45566
    40 = 45567
    41 = 45568
    42 = 45569
    43 = 45570
    44 = 45571
    45 = 45572
    46 = 45573
    47 = 45574
    48 = 45575
}

# Rows whose C/G/J triple uses the "low" quantity figures instead of the
# regular ones -- always the two rows holding the two oldest (smallest)
# date serials in the refreshed window.
$sortedDates = $dates.Values | Sort-Object
$lowDate1 = $sortedDates[0]
$lowDate2 = $sortedDates[1]

$normalC = 0.00170247
$normalG = 465.80531254
$normalJ = 485.38834923

$lowC = 0.00004012
$lowG = 280.99031254
$lowJ = 1941.48834923

$constB = 116.4121952
$constD = 0.008850780000000001
$constE = 0.06933635
$constF = 12792.90181321
$constH = 0.24
$constI = 1.7904431

for ($r = 2; $r -le 48; $r++) {
    $d = $dates[$r]

    if ($d -eq $lowDate1 -or $d -eq $lowDate2) {
        $c = $lowC
        $g = $lowG
        $j = $lowJ
    } else {
        $c = $normalC
        $g = $normalG
        $j = $normalJ
    }

    $ws.Cells.Item($r, 1).Value2 = $d
    $ws.Cells.Item($r, 2).Value2 = $constB
    $ws.Cells.Item($r, 3).Value2 = $c
    $ws.Cells.Item($r, 4).Value2 = $constD
    $ws.Cells.Item($r, 5).Value2 = $constE
    $ws.Cells.Item($r, 6).Value2 = $constF
    $ws.Cells.Item($r, 7).Value2 = $g
    $ws.Cells.Item($r, 8).Value2 = $constH
    $ws.Cells.Item($r, 9).Value2 = $constI
    $ws.Cells.Item($r, 10).Value2 = $j
}

# New rows 46-48 need column A's date style (style index 2 in the
# original file) copied down from the existing dated column.
$ws.Range("A45").Copy() | Out-Null
$ws.Range("A46:A48").PasteSpecial(-4122) | Out-Null
